$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 481.5
$ws.Range("I42").Value = 134.5
$ws.Range("J42").Value = 1002
$ws.Range("K42").Value = 403.5
$ws.Range("L42").Value = 3006
$ws.Range("M42").Value = -173.5
$ws.Range("N42").Value = -3466

$ws.Range("H69").Value = 3450.625
$ws.Range("I69").Value = 3210
$ws.Range("J69").Value = 3485
$ws.Range("K69").Value = 9630
$ws.Range("L69").Value = 10455
$ws.Range("M69").Value = -8756
$ws.Range("N69").Value = -12203

$ws.Range("H72").Value = 3450.625
$ws.Range("I72").Value = 3210
$ws.Range("J72").Value = 3485
$ws.Range("K72").Value = 28890
$ws.Range("L72").Value = 31365
$ws.Range("M72").Value = -24522
$ws.Range("N72").Value = -40101

$ws.Range("H76").Value = 3779
$ws.Range("I76").Value = 3779
$ws.Range("K76").Value = 3779
$ws.Range("M76").Value = -3464

$ws.Range("H79").Value = 3779
$ws.Range("I79").Value = 3779
$ws.Range("K79").Value = 3779
$ws.Range("M79").Value = -2687

$ws.Range("H86").Value = 2019965.2
$ws.Range("I86").Value = 3229464.5
$ws.Range("J86").Value = 4133.1665
$ws.Range("K86").Value = 3229464.5
$ws.Range("L86").Value = 4133.1665
$ws.Range("M86").Value = -3228341.5
$ws.Range("N86").Value = -6379.1665

$ws.Range("H89").Value = 2019965.2
$ws.Range("I89").Value = 3229464.5
$ws.Range("J89").Value = 4133.1665
$ws.Range("K89").Value = 16147322.5
$ws.Range("L89").Value = 20665.8325
$ws.Range("M89").Value = -16141706.5
$ws.Range("N89").Value = -31897.8325

$ws.Range("H107").Value = 540.9
$ws.Range("I107").Value = 545.44446
$ws.Range("K107").Value = 545.44446
$ws.Range("M107").Value = 1374.55554

$ws.Range("H113").Value = 3286.9714
$ws.Range("I113").Value = 3487.9412
$ws.Range("J113").Value = 3097.1667
$ws.Range("K113").Value = 3487.9412
$ws.Range("L113").Value = 3097.1667
$ws.Range("M113").Value = -233.9412000000002
$ws.Range("N113").Value = -9605.1667

$ws.Range("H138").Value = 6240.2246
$ws.Range("I138").Value = 2489.111
$ws.Range("K138").Value = 7467.333
$ws.Range("M138").Value = -2327.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 816.1667
$ws.Range("I88").Value = 692.1429000000001
$ws.Range("J88").Value = 989.8
$ws.Range("K88").Value = 692.1429000000001
$ws.Range("L88").Value = 989.8
$ws.Range("M88").Value = -286.1429000000001
$ws.Range("N88").Value = -1801.8

$ws.Range("H91").Value = 816.1667
$ws.Range("I91").Value = 692.1429000000001
$ws.Range("J91").Value = 989.8
$ws.Range("K91").Value = 692.1429000000001
$ws.Range("L91").Value = 989.8
$ws.Range("M91").Value = 711.8570999999999
$ws.Range("N91").Value = -3797.8

$ws.Range("H102").Value = 979748.5600000001
$ws.Range("I102").Value = 1142605.8
$ws.Range("K102").Value = 1142605.8
$ws.Range("M102").Value = -1140983.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3266.0356
$ws.Range("I86").Value = 2930.8
$ws.Range("K86").Value = 2930.8
$ws.Range("M86").Value = -1807.8

$ws.Range("H89").Value = 3266.0356
$ws.Range("I89").Value = 2930.8
$ws.Range("K89").Value = 14654
$ws.Range("M89").Value = -9038

$ws.Range("H99").Value = 1159392.1
$ws.Range("I99").Value = 1227415.2
$ws.Range("K99").Value = 1227415.2
$ws.Range("M99").Value = -1225917.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 25273
$ws.Range("I99").Value = 44061.332
$ws.Range("K99").Value = 44061.332
$ws.Range("M99").Value = -42563.332

$ws.Range("H126").Value = 25273
$ws.Range("I126").Value = 44061.332
$ws.Range("K126").Value = 132183.996
$ws.Range("M126").Value = -129713.996

$ws.Range("H132").Value = 2776.5925
$ws.Range("I132").Value = 2825.6924
$ws.Range("K132").Value = 8477.0772
$ws.Range("M132").Value = -5947.0772

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 49555704
$ws.Range("I4").Value = 1312843
$ws.Range("K4").Value = 3938529
$ws.Range("M4").Value = -3938417

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1199959.1
$ws.Range("I80").Value = 2378946.5
$ws.Range("J80").Value = 20971.715
$ws.Range("K80").Value = 2378946.5
$ws.Range("L80").Value = 20971.715
$ws.Range("M80").Value = -2377948.5
$ws.Range("N80").Value = -22967.715

$ws.Range("H83").Value = 1199959.1
$ws.Range("I83").Value = 2378946.5
$ws.Range("J83").Value = 20971.715
$ws.Range("K83").Value = 11894732.5
$ws.Range("L83").Value = 104858.575
$ws.Range("M83").Value = -11889740.5
$ws.Range("N83").Value = -114842.575

$ws.Range("H102").Value = 3416.1936
$ws.Range("I102").Value = 2229.3333
$ws.Range("K102").Value = 2229.3333
$ws.Range("M102").Value = -607.3332999999998

$ws.Range("H122").Value = 848832.1
$ws.Range("I122").Value = 1002756.2
$ws.Range("J122").Value = 2250
$ws.Range("K122").Value = 3008268.6
$ws.Range("L122").Value = 6750
$ws.Range("M122").Value = -3005818.6
$ws.Range("N122").Value = -11650

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6697.6665
$ws.Range("I46").Value = 4741.6665
$ws.Range("K46").Value = 4741.6665
$ws.Range("M46").Value = -4553.6665

$ws.Range("H68").Value = 735202.1
$ws.Range("I68").Value = 1516892.5
$ws.Range("J68").Value = 2367.4375
$ws.Range("K68").Value = 1516892.5
$ws.Range("L68").Value = 2367.4375
$ws.Range("M68").Value = -1516143.5
$ws.Range("N68").Value = -3865.4375

$ws.Range("H71").Value = 735202.1
$ws.Range("I71").Value = 1516892.5
$ws.Range("J71").Value = 2367.4375
$ws.Range("K71").Value = 7584462.5
$ws.Range("L71").Value = 11837.1875
$ws.Range("M71").Value = -7580718.5
$ws.Range("N71").Value = -19325.1875

$ws.Range("H93").Value = 1490.0834
$ws.Range("I93").Value = 1547.625
$ws.Range("K93").Value = 1547.625
$ws.Range("M93").Value = -299.625

$ws.Range("H100").Value = 1848.2858
$ws.Range("I100").Value = 1823.1666
$ws.Range("K100").Value = 1823.1666
$ws.Range("M100").Value = -1282.1666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3792135.8
$ws.Range("J81").Value = 6946555
$ws.Range("L81").Value = 13893110
$ws.Range("N81").Value = -13895232

$ws.Range("H84").Value = 3792135.8
$ws.Range("J84").Value = 6946555
$ws.Range("L84").Value = 69465550
$ws.Range("N84").Value = -69476158

$ws.Range("H126").Value = 1307.7222
$ws.Range("I126").Value = 1250
$ws.Range("J126").Value = 1398.4286
$ws.Range("K126").Value = 3750
$ws.Range("L126").Value = 4195.2858
$ws.Range("M126").Value = -1280
$ws.Range("N126").Value = -9135.2858

$ws.Range("H136").Value = 7429
$ws.Range("I136").Value = 2640.625
$ws.Range("J136").Value = 8941.118
$ws.Range("K136").Value = 7921.875
$ws.Range("L136").Value = 26823.354
$ws.Range("M136").Value = -5371.875
$ws.Range("N136").Value = -31923.354
